$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Juan Pablo", "Bautista Clavijo", "sslsc@udistrital.edu.co", "1730262733833-reporte_estudiante (2).pdf"),
    @("Juan Pablo", "Bautista Clavijo", "sslsc@udistrital.edu.co", "1730262850260-Escudo_UD.png"),
    @("Juan Pablo", "Bautista Clavijo", "jbautistaclavijo@gmail.com", "1730262972598-reporte_estudiante (2).pdf"),
    @("Juan Pablo", "Bautista Clavijo", "jbautistaclavijo@gmail.com", "1730263075607-reporte_estudiante (2).pdf"),
    @("Juan Pablo", "Bautista Clavijo", "jbautistaclavijo@gmail.com", "1730263362238-reporte_estudiante (2).pdf"),
    @("Juan Pablo", "Bautista Clavijo", "jbautistaclavijo@gmail.com", "1730263606699-Ejercicios ED HomogeÌneas.pdf"),
    @("Juan Pablo", "Bautista Clavijo", "jbautistaclavijo@gmail.com", "1730263744368-reporte_estudiante (2).pdf"),
    @("Juan Pablo", "Bautista Clavijo", "jbautistaclavijo@gmail.com", "1730263775745-reporte_estudiante (2).pdf"),
    @("Juan Pablo", "Bautista Clavijo", "jbautistaclavijo@gmail.com", "1730264179945-reporte_estudiante (2).pdf"),
    @("Juan Pablo", "Bautista Clavijo", "jbautistaclavijo@gmail.com", "1730264270327-reporte_estudiante (2).pdf"),
    @("Juan Pablo", "Bautista Clavijo", "jbautistaclavijo@gmail.com", "1730265483895-reporte_estudiante (2).pdf"),
    @("Juan Pablo", "Bautista Clavijo", "jbautistaclavijo@gmail.com", "1730265556103-Ejercicios ED Exactas.pdf"),
    @("Juan Pablo", "Bautista Clavijo", "jbautistaclavijo@gmail.com", "1730265720556-reporte_estudiante (2).pdf"),
    @("Juan Pablo", "Bautista Clavijo", "jbautistaclavijo@gmail.com", "1730265909314-reporte_estudiante (2).pdf"),
    @("Juan Pablo", "Bautista Clavijo", "jbautistaclavijo@gmail.com", "1730266029793-reporte_estudiante (2).pdf"),
    @("Juan Pablo", "Bautista Clavijo", "jbautistaclavijo@gmail.com", "1730266224051-reporte_estudiante (2).pdf")
)

$startRow = 4
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
